$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.237585067749023
$ws.Range("B1").Value = 2.199332237243652
$ws.Range("C1").Value = 6.023508071899414
$ws.Range("D1").Value = 1.251180171966553
$ws.Range("E1").Value = 1.30712902545929
